$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test cases")
$ws.Activate()

# ------------------------------------------------------------------
# New test cases appended to the "Test cases" sheet:
#   rows 94-97  -> UC21 / Function 21: Add comment
#   rows 98-101 -> UC22 / Function 22: Add rating
#
# Cell values are written in the same order the original author must
# have typed them so the shared-string table ends up in the same
# order (UC21, "Function 21: Add comment", UC22,
# "Function 22: Add rating").
# ------------------------------------------------------------------

# --- Row 95: first sub-row of the UC21 group (plain cells, no
#     pre-existing row formatting) -----------------------------------
$ws.Range("A95").Value2 = "UC21"
$ws.Range("B95").Value = "'UI01"

# --- Row 94: the "Function 21" section header row, formatted like
#     the previous section header (row 90) ---------------------------
$ws.Range("A90:J90").Copy()
$ws.Range("A94").PasteSpecial(-4122)
$ws.Range("A94").Value2 = "Function 21: Add comment"

# --- Rows 96-97: remaining UC21 sub-rows -----------------------------
$ws.Range("A96").Value2 = "UC21"
$ws.Range("B96").Value2 = "UI02"

$ws.Range("A97").Value2 = "UC21"
$ws.Range("B97").Value2 = "UI03"

# --- Row 99: first sub-row of the UC22 group -------------------------
$ws.Range("A99").Value2 = "UC22"
$ws.Range("B99").Value = "'UI01"

# --- Row 98: the "Function 22" section header row --------------------
$ws.Range("A90:J90").Copy()
$ws.Range("A98").PasteSpecial(-4122)
$ws.Range("A98").Value2 = "Function 22: Add rating"

# --- Rows 100-101: remaining UC22 sub-rows ---------------------------
$ws.Range("A100").Value2 = "UC22"
$ws.Range("B100").Value2 = "UI02"

$ws.Range("A101").Value2 = "UC22"
$ws.Range("B101").Value2 = "UI03"

$excel.CutCopyMode = $false

$ws.Range("D101").Select() | Out-Null
